$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.083.83'
$ws.Range("E2").Value = '  -4.21%  '

$ws.Range("D3").Value = '3.004.14'
$ws.Range("E3").Value = '  -6.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.72'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("D8").Value = '2.989.91'
$ws.Range("E8").Value = '  -6.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.473'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -12.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -12.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.443'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -10.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -13.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000210'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -13.70%  '

$ws.Range("D15").Value = '3.485.80'
$ws.Range("E15").Value = '  -6.63%  '

$ws.Range("D16").Value = '64.115.60'
$ws.Range("E16").Value = '  -4.26%  '

$ws.Range("E17").Value = '  -4.31%  '

$ws.Range("D18").Value = '3.003.64'
$ws.Range("E18").Value = '  -6.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '477.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -11.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -11.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.648'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -14.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -14.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -10.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.77%  '

$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("E27").Value = '  -15.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -13.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.997'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '504.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -11.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -12.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0408'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.120'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0764'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -13.93%  '

$ws.Range("D42").Value = '2.759.46'
$ws.Range("E42").Value = '  -5.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.88%  '

$ws.Range("E45").Value = '  -11.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.73%  '

$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").Value = '0.0₃0501'
$ws.Range("E47").Value = '  -13.43%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '115.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.49%  '

$ws.Range("E49").Value = '  -9.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -14.02%  '

$ws.Range("E51").Value = '  -18.32%  '
